$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate header row (row 1) labels to Chinese
$ws.Range("A1").Value = "證券編碼"
$ws.Range("B1").Value = "最新價格"
$ws.Range("C1").Value = "最新價格時間戳"
$ws.Range("D1").Value = "今日漲跌價格"
$ws.Range("E1").Value = "今日漲跌幅度"
$ws.Range("F1").Value = "前日收盤價"
$ws.Range("G1").Value = "市場識別碼"
$ws.Range("I1").Value = "最高價"
$ws.Range("J1").Value = "最低價"
$ws.Range("K1").Value = "價格定點"
$ws.Range("L1").Value = "以百分比交易"
$ws.Range("M1").Value = "交易結束時間"
$ws.Range("N1").Value = "交易開始時間"
$ws.Range("O1").Value = "成交額（歐元）"
$ws.Range("P1").Value = "成交量（件數）"
$ws.Range("Q1").Value = "名義成交額"
$ws.Range("R1").Value = "52週最高價"
$ws.Range("S1").Value = "52週最低價"
$ws.Range("T1").Value = "貨幣"
$ws.Range("U1").Value = "最小可交易單位"

# Update data row (row 2) values
$ws.Range("B2").Value = 91
$ws.Range("C2").Value = "2025-02-21 16:37"
$ws.Range("D2").Value = 0.5600000000000001
$ws.Range("E2").Value = 0.62
$ws.Range("I2").Value = 91
$ws.Range("J2").Value = 91
